# Auto-generated update of Moogle Profits market data cells.
# Applies the per-cell value changes described by the upstream XML diff
# (scheduled runner refresh of currentAveragePrice / LevePrice / LeveProfit columns).
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 38 (source diff hunk starting at line 2509)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4272.88
$ws.Range("J38").Value = 7586.9165
$ws.Range("L38").Value = 22760.7495
$ws.Range("N38").Value = -23504.7495

# Sheet ALC, row 100 (source diff hunk starting at line 5601)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4049.7273
$ws.Range("I100").Value = 3881
$ws.Range("K100").Value = 3881
$ws.Range("M100").Value = -3340

# Sheet ALC, row 116 (source diff hunk starting at line 6388)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 18809.637
$ws.Range("J116").Value = 19127.666
$ws.Range("L116").Value = 19127.666
$ws.Range("N116").Value = -26011.666

# Sheet ALC, row 125 (source diff hunk starting at line 6823)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 73308.86
$ws.Range("I125").Value = 2053.9
$ws.Range("J125").Value = 251446.25
$ws.Range("K125").Value = 18485.1
$ws.Range("L125").Value = 2263016.25
$ws.Range("M125").Value = -16025.1
$ws.Range("N125").Value = -2267936.25

# Sheet ALC, row 132 (source diff hunk starting at line 7166)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1825.1273
$ws.Range("I132").Value = 1737.4902
$ws.Range("K132").Value = 5212.4706
$ws.Range("M132").Value = -2682.4706

# Sheet ALC, row 137 (source diff hunk starting at line 7417)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2176.276
$ws.Range("I137").Value = 2081.2307
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 6243.6921
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -3693.6921
$ws.Range("N137").Value = -14100

# Sheet ALC, row 138 (source diff hunk starting at line 7469)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3961.6575
$ws.Range("I138").Value = 3374.152
$ws.Range("J138").Value = 4962.593
$ws.Range("K138").Value = 10122.456
$ws.Range("L138").Value = 14887.779
$ws.Range("M138").Value = -4982.456
$ws.Range("N138").Value = -25167.779

# Sheet ARM, row 32 (source diff hunk starting at line 9217)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7026.742
$ws.Range("I32").Value = 3069.291
$ws.Range("K32").Value = 3069.291
$ws.Range("M32").Value = -2782.291

# Sheet ARM, row 101 (source diff hunk starting at line 12544)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 39950
$ws.Range("J101").Value = 39950
$ws.Range("L101").Value = 39950
$ws.Range("N101").Value = -46440

# Sheet ARM, row 102 (source diff hunk starting at line 12593)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1525
$ws.Range("I102").Value = 1525
$ws.Range("K102").Value = 1525
$ws.Range("M102").Value = 97

# Sheet BSM, row 82 (source diff hunk starting at line 18462)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 48767.727
$ws.Range("I82").Value = 13778.571
$ws.Range("J82").Value = 109998.75
$ws.Range("K82").Value = 13778.571
$ws.Range("L82").Value = 109998.75
$ws.Range("M82").Value = -13395.571
$ws.Range("N82").Value = -110764.75

# Sheet BSM, row 85 (source diff hunk starting at line 18615)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 48767.727
$ws.Range("I85").Value = 13778.571
$ws.Range("J85").Value = 109998.75
$ws.Range("K85").Value = 13778.571
$ws.Range("L85").Value = 109998.75
$ws.Range("M85").Value = -12452.571
$ws.Range("N85").Value = -112650.75

# Sheet BSM, row 86 (source diff hunk starting at line 18667)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3217.5
$ws.Range("I86").Value = 3638.1667
$ws.Range("J86").Value = 1955.5
$ws.Range("K86").Value = 3638.1667
$ws.Range("L86").Value = 1955.5
$ws.Range("M86").Value = -2515.1667
$ws.Range("N86").Value = -4201.5

# Sheet BSM, row 89 (source diff hunk starting at line 18814)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3217.5
$ws.Range("I89").Value = 3638.1667
$ws.Range("J89").Value = 1955.5
$ws.Range("K89").Value = 18190.8335
$ws.Range("L89").Value = 9777.5
$ws.Range("M89").Value = -12574.8335
$ws.Range("N89").Value = -21009.5

# Sheet BSM, row 134 (source diff hunk starting at line 20986)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2940.2144
$ws.Range("I134").Value = 1512.4166
$ws.Range("J134").Value = 11507
$ws.Range("K134").Value = 4537.2498
$ws.Range("L134").Value = 34521
$ws.Range("M134").Value = -2002.2498
$ws.Range("N134").Value = -39591

# Sheet CRP, row 94 (source diff hunk starting at line 26004)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5702.364
$ws.Range("I94").Value = 5857.6665
$ws.Range("J94").Value = 5516
$ws.Range("K94").Value = 5857.6665
$ws.Range("L94").Value = 5516
$ws.Range("M94").Value = -5406.6665
$ws.Range("N94").Value = -6418

# Sheet CRP, row 96 (source diff hunk starting at line 26102)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 3837.75
$ws.Range("J96").Value = 3837.75
$ws.Range("L96").Value = 3837.75
$ws.Range("N96").Value = -9329.75

# Sheet CRP, row 122 (source diff hunk starting at line 27349)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3440.889
$ws.Range("I122").Value = 3225.32
$ws.Range("K122").Value = 9675.960000000001
$ws.Range("M122").Value = -7225.960000000001

# Sheet CRP, row 132 (source diff hunk starting at line 27842)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5728.3887
$ws.Range("I132").Value = 3392.2307
$ws.Range("K132").Value = 10176.6921
$ws.Range("M132").Value = -7646.6921

# Sheet CRP, row 133 (source diff hunk starting at line 27894)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# Sheet CUL, row 2 (source diff hunk starting at line 28435)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 997.1539
$ws.Range("I2").Value = 160.11111
$ws.Range("J2").Value = 2880.5
$ws.Range("K2").Value = 960.66666
$ws.Range("L2").Value = 17283
$ws.Range("M2").Value = -847.66666
$ws.Range("N2").Value = -17509

# Sheet CUL, row 76 (source diff hunk starting at line 32160)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 14034.077
$ws.Range("J76").Value = 16810.834
$ws.Range("L76").Value = 50432.50199999999
$ws.Range("N76").Value = -51198.50199999999

# Sheet CUL, row 79 (source diff hunk starting at line 32316)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 14034.077
$ws.Range("J79").Value = 16810.834
$ws.Range("L79").Value = 50432.50199999999
$ws.Range("N79").Value = -53084.50199999999

# Sheet CUL, row 113 (source diff hunk starting at line 33991)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4348702
$ws.Range("J113").Value = 4831836
$ws.Range("L113").Value = 14495508
$ws.Range("N113").Value = -14499848

# Sheet GSM, row 49 (source diff hunk starting at line 37830)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 41900
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 41900
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 41900
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -42268

# Sheet GSM, row 97 (source diff hunk starting at line 40125)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 869.8333
$ws.Range("I97").Value = 803
$ws.Range("K97").Value = 803
$ws.Range("M97").Value = -307

# Sheet GSM, row 102 (source diff hunk starting at line 40376)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3140.697
$ws.Range("I102").Value = 1817.4762
$ws.Range("K102").Value = 1817.4762
$ws.Range("M102").Value = -195.4762000000001

# Sheet GSM, row 122 (source diff hunk starting at line 41329)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3144.88
$ws.Range("I122").Value = 1559.2632
$ws.Range("J122").Value = 8166
$ws.Range("K122").Value = 4677.7896
$ws.Range("L122").Value = 24498
$ws.Range("M122").Value = -2227.7896
$ws.Range("N122").Value = -29398

# Sheet LTW, row 2 (source diff hunk starting at line 42391)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 37999.75
$ws.Range("J2").Value = 37999.75
$ws.Range("L2").Value = 37999.75
$ws.Range("N2").Value = -38223.75

# Sheet LTW, row 40 (source diff hunk starting at line 44250)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7247.684
$ws.Range("I40").Value = 4041.8823
$ws.Range("K40").Value = 4041.8823
$ws.Range("M40").Value = -3905.8823

# Sheet LTW, row 100 (source diff hunk starting at line 47139)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

# Sheet LTW, row 132 (source diff hunk starting at line 48686)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1993.1915
$ws.Range("I132").Value = 1255.4166
$ws.Range("J132").Value = 4407.727
$ws.Range("K132").Value = 3766.2498
$ws.Range("L132").Value = 13223.181
$ws.Range("M132").Value = -1236.2498
$ws.Range("N132").Value = -18283.181

# Sheet LTW, row 136 (source diff hunk starting at line 48885)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5927.9644
$ws.Range("I136").Value = 3198.6191
$ws.Range("K136").Value = 9595.8573
$ws.Range("M136").Value = -7045.8573

# Sheet WVR, row 122 (source diff hunk starting at line 55117)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2130.389
$ws.Range("I122").Value = 2476.3333
$ws.Range("J122").Value = 1438.5
$ws.Range("K122").Value = 7428.999899999999
$ws.Range("L122").Value = 4315.5
$ws.Range("M122").Value = -4978.999899999999
$ws.Range("N122").Value = -9215.5

# Sheet WVR, row 132 (source diff hunk starting at line 55604)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10239.875
$ws.Range("I132").Value = 5584.2
$ws.Range("J132").Value = 17999.334
$ws.Range("K132").Value = 16752.6
$ws.Range("L132").Value = 53998.00199999999
$ws.Range("M132").Value = -14222.6
$ws.Range("N132").Value = -59058.00199999999
